# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# The "K" column (column G) holds a recalculated statistic for each trade
# row. This script rewrites the previously-stored "Strike#"-derived values
# in column G with the freshly calculated K values, row by row, leaving
# every other column untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new K value
$kValues = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 1
    6  = 1
    8  = 2
    9  = 0
    10 = 1
    11 = 1
    12 = 1
    13 = 2
    14 = 0
    15 = 0
    16 = 2
    17 = 1
    18 = 2
    19 = 1
    20 = 2
    21 = 2
    22 = 2
    23 = 3
    24 = 1
    25 = 1
    26 = 1
    27 = 1
    28 = 0
    29 = 2
    30 = 0
    31 = 0
    32 = 1
    33 = 2
    34 = 1
    35 = 3
    36 = 2
    37 = 4
    38 = 1
    39 = 1
    40 = 2
    41 = 2
    42 = 1
    43 = 1
    44 = 0
    45 = 1
    46 = 0
    47 = 0
    48 = 0
    49 = 0
    50 = 1
    51 = 0
    52 = 1
    53 = 1
    54 = 2
    55 = 0
    56 = 1
    57 = 0
    58 = 1
    59 = 0
    60 = 1
    61 = 0
    62 = 1
    63 = 0
    64 = 1
    66 = 1
    68 = 2
    69 = 2
    70 = 2
    71 = 0
    72 = 3
    73 = 0
    74 = 1
    75 = 1
    76 = 3
    77 = 1
    78 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
